$d = $word.ActiveDocument

# The document currently ends with an empty paragraph immediately before
# the final sectPr. We need to keep that empty paragraph intact and add
# the new "Conclusions and Future Works" section (heading + two
# sub-sections with body text) right after it.

$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newLastPara = $d.Paragraphs.Last
$insertRange = $newLastPara.Range

$fragment = @'
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
    </w:rPr>
    <w:t>CONCLUSIONS AND FUTURE WORKS:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>a. Suggested Improvements:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">It is imperative to incorporate additional dynamic aspects into the process </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>in order to</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> envisage future advances. The accuracy of workload assessments could be improved by incorporating machine learning algorithms to adjust to shifting work patterns, personal preferences, and organizational dynamics. Investigating technological solutions that promote smooth communication and cooperation may also help to improve project management effectiveness. Customized time management tools and ongoing feedback loops for </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>reevaluation</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> may provide professionals specialized solutions. Moreover, a more thorough examination of factors related to psychological well-being in the context of managing workloads may provide a comprehensive viewpoint.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>b. Limitations to Solution:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">   </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Although the suggested solution offers a strong structure, its application can be limited in some situations. The methodology may face difficulties in extremely unexpected and quickly evolving work situations where tasks are intrinsically dynamic and cannot be </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>planned ahead</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> of time. Furthermore, it could be difficult for companies with strict deadlines and regulatory frameworks to adopt the flexibility recommended for retaining free time. It may be difficult for people with little autonomy or decision-making authority to put delegation techniques into practice. In situations when technology tools are not widely accepted or available, the solution could potentially have drawbacks.</w:t>
  </w:r>
</w:p>
'@

$insertRange.InsertXML($fragment)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
